$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new row at position 109. This shifts the (previously blank) rows
# 109-111 down to 110-112 and the summary rows 112-113 down to 113-114,
# with relative formula references auto-adjusting (e.g. B112 -> B113 inside
# the old B113 formula).
# ---------------------------------------------------------------------------
$ws.Rows.Item(109).Insert()

# ---------------------------------------------------------------------------
# Row 109 is the brand-new "latest" row, so it should inherit the highlight
# ("Good"/green) formatting that row 108 currently has, before row 108 gets
# reformatted down to the regular ("Neutral"/yellow) look.
# ---------------------------------------------------------------------------
$ws.Range("A108:H108").Copy()
$ws.Range("A109:H109").PasteSpecial(-4122)
$ws.Range("I108").Copy()
$ws.Range("I109").PasteSpecial(-4122)

# Row 108 now becomes a "past" row and should look like row 107 (Neutral).
$ws.Range("A107:H107").Copy()
$ws.Range("A108:H108").PasteSpecial(-4122)

# The "daily rate to achieve June 20 target" note moves from row 108 to the
# new latest row, 109.
$ws.Range("I108").Clear()

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Fill in the new day's data on row 109.
# ---------------------------------------------------------------------------
$ws.Range("A109").Value = 44294
$ws.Range("B109").Value = 2524
$ws.Range("C109").Formula = "=(AVERAGE(B103:B109))"
$ws.Range("D109").Formula = "=(D108-B109)"
$ws.Range("E109").Formula = "=E108+B109"
$ws.Range("F109").Formula = "=D109/C109"
$ws.Range("G109").Formula = "=A109+F109"
$ws.Range("H109").Formula = "=D109/84"
$ws.Range("I109").Value = "daily rate to achieve June 20 target"

# ---------------------------------------------------------------------------
# Update the view: scrolled down a bit further and a new active selection.
# ---------------------------------------------------------------------------
$win = $excel.Windows.Item(1)
$win.ScrollRow = 85
$ws.Range("I113").Select()
